# jiga7 | Mise en page
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Colour the section heading "Paramétrage de l'application" in red
#    (applies both to the run text and to the paragraph mark, i.e.
#    the <w:color w:val="FF0000"/> added to both rPr blocks).
# ------------------------------------------------------------------
$headingPara = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*Paramétrage de l’application*") {
    $headingPara = $p
  }
}
if ($headingPara -ne $null) {
  $headingPara.Range.Font.Color = 255
}

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document to the
#    point right after "...cliquer sur" (before the trailing space),
#    in the "le système affiche ... cliquer sur Dépôt" bullet. This
#    also naturally splits that run into "...cliquer sur" + " ".
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("cliquer sur", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterSur = $find.Parent.End

$bookmarkRange = $d.Range($afterSur, $afterSur)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Re-touch the "...cliquer sur" run so it is re-serialized without a
# spurious xml:space="preserve" now that it no longer has a trailing
# space (the bookmark insertion above would otherwise leave a stray
# xml:space="preserve" on that run).
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute("cliquer sur", $true, $false, $false, $false, $false, $true, 1, $false, "cliquer sur", 2) | Out-Null
